# Applies:
#  1. Yellow-highlight the "Asegurarse que se puede organizar..." bullet paragraph.
#  2. Yellow-highlight the "El tamaño de la fuente de la caja de texto..." bullet paragraph.
#  3. Merge the split comment run ("Mucho mas complicado..." + "DialogChild") back into one run
#     (drops the spell-check markers around "DialogChild").

$d = $word.ActiveDocument

# --- 1 & 2: highlight whole paragraphs (all runs) in yellow ---------------

function Highlight-ParagraphByText($needle) {
    $rng = $d.Content
    $found = $rng.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1).Range
        $para.HighlightColorIndex = 7  # wdYellow
    }
}

Highlight-ParagraphByText("Asegurarse que se puede organizar")
Highlight-ParagraphByText("El tamaño de la fuente de la caja de texto")

# --- 3: merge the comment's split run back together ------------------------

$comments = $d.Comments
for ($i = 1; $i -le $comments.Count; $i++) {
    $c = $comments.Item($i)
    $cr = $c.Range
    if ($cr.Text -like "*DialogChild*") {
        $cr.Find.Execute("Mucho mas complicado de resolver, dado que la validación se hace al cerrar el DialogChild", $true, $false, $false, $false, $false, $true, 1, $false, "Mucho mas complicado de resolver, dado que la validación se hace al cerrar el DialogChild", 2)
    }
}
